$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1.xml): update "want to go" counts (column F)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 19
$wsExhibit.Range("F4").Value = 3347
$wsExhibit.Range("F6").Value = 30
$wsExhibit.Range("F7").Value = 168

# "全部类型" sheet (sheet4.xml): same events repeated, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 19
$wsAll.Range("F8").Value = 3347
$wsAll.Range("F10").Value = 30
$wsAll.Range("F12").Value = 168
